$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# Change 1: the empty paragraph right after the intro text used to hold only
# the (hidden, system) "_GoBack" bookmark. Remove the bookmark so the
# paragraph becomes genuinely empty.
# ---------------------------------------------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# ---------------------------------------------------------------------------
# Change 2: in the "Sessione" table, the "Ruolo (chiave primaria di
# utenteModeratore)" field name is renamed to "Username", keeping the
# parenthetical note, now split across two runs.
# ---------------------------------------------------------------------------
$tSessione = $d.Tables.Item(7)
$cellRuolo = $tSessione.Cell(5, 1)
$pRuolo = $cellRuolo.Range.Paragraphs.Item(1)
$xmlRuolo = "<w:p $wNs><w:r><w:t>Username</w:t></w:r><w:r><w:t xml:space=`"preserve`"> (chiave primaria di utenteModeratore)</w:t></w:r></w:p>"
$pRuolo.Range.InsertXML($xmlRuolo)

# ---------------------------------------------------------------------------
# Change 3: in the "Realizza (relazione tra utenteModeratore e Storia)"
# table, the "Ruolo" field is replaced by a "Username" field that is the
# primary key of utenteRegistrato (with the _GoBack bookmark moved here),
# and its constraints cell grows an extra paragraph describing the new
# field.
# ---------------------------------------------------------------------------
$tRealizza = $d.Tables.Item(12)

$cellUsername = $tRealizza.Cell(3, 1)
$pUsername = $cellUsername.Range.Paragraphs.Item(1)
$xmlUsername = "<w:p $wNs><w:r><w:t>Username</w:t></w:r><w:r><w:t xml:space=`"preserve`"> </w:t></w:r><w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/><w:bookmarkEnd w:id=`"0`"/><w:r><w:t xml:space=`"preserve`">(chiave primaria di </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>utenteRegistrato</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t>)</w:t></w:r></w:p>"
$pUsername.Range.InsertXML($xmlUsername)

$tRealizza2 = $d.Tables.Item(12)
$cellVincoli = $tRealizza2.Cell(3, 2)
$pVincoli = $cellVincoli.Range.Paragraphs.Item(1)
$xmlVincoli = "<w:p $wNs><w:r><w:t xml:space=`"preserve`">Lunghezza massima: 15 caratteri, </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>unique</w:t></w:r><w:proofErr w:type=`"spellEnd`"/></w:p><w:p $wNs><w:r><w:t xml:space=`"preserve`">Chiave </w:t></w:r><w:r><w:t>esterna</w:t></w:r></w:p>"
$pVincoli.Range.InsertXML($xmlVincoli)
